$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row with the second modification note
$ws.Range("A3").Value = "Segunda modificacion para fines practicos. 10/19/18 8:05"

# Move the active selection to A4, matching the post-edit state
$ws.Range("A4").Select()
